$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New insanity card: Ophidiophobia (row 26)
$ws.Range("A26").Value = "Ophidiophobia"
$ws.Range("B26").Value = 3
$ws.Range("C26").Value = "Why did it have to be snakes?"
$ws.Range("D26").Value = "<p>You cannot voluntarily move into a space that contains a Hunting Horror, Feathered Serpent, Serpent Person, or Ancient Basilisk. You win or lose the game as normal.</p>"
$ws.Rows.Item(26).RowHeight = 45

# New insanity card: Treasure Hunter (row 27)
$ws.Range("A27").Value = "Treasure Hunter"
$ws.Range("B27").Value = 2
$ws.Range("C27").Value = "You deserve a reward for your good work. There must be some valuables around here somewhere."
$ws.Range("D27").Value = "<p>You do not win the game as normal. Instead, you win only if the investigation is complete and there are no Items on the board.</p>"
$ws.Rows.Item(27).RowHeight = 30

# Match the saved view: selection moved to the new first empty row
$ws.Range("A28").Select()
